$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.061818241674302499
$ws.Range("A2").Value = -0.036376383531251832
$ws.Range("A3").Value = -0.0089999994876368561
$ws.Range("A4").Value = 0.28399094447783924
$ws.Range("A5").Value = -0.0059999995038975129
$ws.Range("A6").Value = -0.037420417724664645
$ws.Range("A7").Value = -0.019999999397656509
$ws.Range("A8").Value = 0.019302722220913893
$ws.Range("A9").Value = -0.0059999994838229043
$ws.Range("A10").Value = -0.0059999994815740365
$ws.Range("A11").Value = -0.0044999994911414376
$ws.Range("A12").Value = -0.0059999994812791613
$ws.Range("A13").Value = -0.0059999994808510593
$ws.Range("A14").Value = -0.011999999442631193
$ws.Range("A15").Value = -0.005999999481974605
$ws.Range("A16").Value = 0.036816884431109642
$ws.Range("A17").Value = -0.0059999994815846946
$ws.Range("A18").Value = -0.0089999994618628065
$ws.Range("A19").Value = -0.0089999994875329392
$ws.Range("A20").Value = -0.0089999994830627372
$ws.Range("A21").Value = -0.0089999994824063734
$ws.Range("A22").Value = -0.0089999994819125462
$ws.Range("A23").Value = -0.0089999994829232932
$ws.Range("A24").Value = -0.041999999261983056
$ws.Range("A25").Value = -0.043831854645389079
$ws.Range("A26").Value = -0.0059999994857449224
$ws.Range("A27").Value = -0.0059999994830617354
$ws.Range("A28").Value = -0.0059999994717365723
$ws.Range("A29").Value = -0.011999999425324148
$ws.Range("A30").Value = -0.019999999370311272
$ws.Range("A31").Value = -0.014999999396946961
$ws.Range("A32").Value = -0.020999999357614207
$ws.Range("A33").Value = -0.0059999994534685186
